{"js": "// Remove the last table row (\"Regras de neg\u00f3cio\" / \"Regras de neg\u00f3cio do\n// cliente situadas no artefato 19\") from the references table, per the\n// commit \"corre\u00e7\u00e3o para nota da prova\".\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  // Walk rows back-to-front so deleting doesn't disturb earlier indices,\n  // and match on the row's own text so we only touch the intended row.\n  for (let i = rows.items.length - 1; i >= 0; i--) {\n    const row = rows.items[i];\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    let rowText = \"\";\n    for (const cell of cells.items) {\n      cell.body.load(\"text\");\n    }\n    await context.sync();\n    for (const cell of cells.items) {\n      rowText += cell.body.text;\n    }\n\n    if (\n      rowText.indexOf(\"Regras de neg\u00f3cio\") !== -1 &&\n      rowText.indexOf(\"Regras de neg\u00f3cio do cliente situadas no artefato 19\") !== -1\n    ) {\n      row.delete();\n      await context.sync();\n    }\n  }\n}\n", "ps1": "# Remove the last table row (\"Regras de neg\u00f3cio\" / \"Regras de neg\u00f3cio do\n# cliente situadas no artefato 19\") from the references table, per the\n# commit \"corre\u00e7\u00e3o para nota da prova\".\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    # Walk rows back-to-front so deleting a row doesn't shift the indices\n    # of the rows still to be inspected.\n    for ($i = $t.Rows.Count; $i -ge 1; $i--) {\n        $row = $t.Rows.Item($i)\n        $label = $row.Cells.Item(1).Range.Text\n        $value = $row.Cells.Item(2).Range.Text\n\n        if (($label -like \"*Regras de neg\u00f3cio*\") -and ($value -like \"*Regras de neg\u00f3cio do cliente situadas no artefato 19*\")) {\n            $row.Delete()\n        }\n    }\n}\n"}
